$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experiment_specification")

$rows = @(29, 30, 32, 33, 35, 36)
foreach ($r in $rows) {
    $dCell = $ws.Range("D$r")
    $eCell = $ws.Range("E$r")
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}

[void]$ws.Activate()
[void]$ws.Range("D35:D36").Select()
Write-Output "done"
